# 自动更新Excel文件 - 每日剩余天数递减，到期后重置周期
# For every data row: if remaining days (E) is down to 1, the stay cycle is
# considered finished today, so it resets to the full duration (D) and the
# start date (F) rolls to "today" (2025-12-02). Otherwise E just decrements
# by one day. Rows whose start date isn't a well-formed YYYYMMDD value are
# left untouched (can't compute a cycle for them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251202

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    $totalDays = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startDate = $ws.Cells.Item($r, 6).Value2

    if ($null -eq $totalDays -or $null -eq $remaining -or $null -eq $startDate) {
        continue
    }

    # Validate F looks like an 8-digit YYYYMMDD date; skip malformed rows
    # (e.g. a stray extra digit) exactly like the upstream updater does.
    $dateText = [string][int]$startDate
    if ($dateText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value = $totalDays
        $ws.Cells.Item($r, 6).Value = $today
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
